$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.099494457244873
$ws.Range("B1").Value = 2.123276472091675
$ws.Range("C1").Value = 9.267221450805664
$ws.Range("D1").Value = 2.412626504898071
$ws.Range("E1").Value = 1.295566916465759
